$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column widths: A wider to fit long institution names ---
$ws.Columns.Item(1).ColumnWidth = 50.0

# --- Clear stale trailing data from rows 53-54 (their content moved up into 51-52) ---
$ws.Range("A53").Clear()
$ws.Range("A54").Clear()

# --- Header row: rename A1, add "replace" header in B1 matching A1 (bold Arial) styling ---
$ws.Range("A1").Value = "replacement"
$ws.Range("B1").Value = "replace"
$ws.Range("B1").Font.Bold = $true
$ws.Range("B1").Font.Name = "Arial"

# --- Set cell text values per new layout (rows 2-52) ---
$ws.Range("A1").Value = "replacement"
$ws.Range("B1").Value = "replace"
$ws.Range("A2").Value = "California Institute of Technology"
$ws.Range("A3").Value = "Colorado State University"
$ws.Range("A4").Value = "CUNY"
$ws.Range("A5").Value = "Georgia Institute of Technology"
$ws.Range("A6").Value = "Indiana University"
$ws.Range("A7").Value = "Kent State University"
$ws.Range("A8").Value = "New Mexico State University"
$ws.Range("A9").Value = "North Dakota State University"
$ws.Range("A10").Value = "OSU"
$ws.Range("B10").Value = "Ohio State University-Main Campus"
$ws.Range("A11").Value = "Ohio University"
$ws.Range("A12").Value = "Purdue University"
$ws.Range("A13").Value = "Rutgers University"
$ws.Range("A14").Value = "Southern Illinois University"
$ws.Range("A15").Value = "TX A&M"
$ws.Range("B15").Value = "Texas A & M University"
$ws.Range("A16").Value = "The University of Tennessee"
$ws.Range("A17").Value = "U TX Austin"
$ws.Range("B17").Value = "The University of Texas"
$ws.Range("A18").Value = "UCB"
$ws.Range("B18").Value = "University of California-Berkeley"
$ws.Range("A19").Value = "UCD"
$ws.Range("B19").Value = "University of California-Davis"
$ws.Range("A20").Value = "UCI"
$ws.Range("B20").Value = "University of California-Irvine"
$ws.Range("A21").Value = "UCLA"
$ws.Range("B21").Value = "University of California-Los Angeles"
$ws.Range("A22").Value = "UCM"
$ws.Range("B22").Value = "University of California-Merced"
$ws.Range("A23").Value = "UCR"
$ws.Range("B23").Value = "University of California-Riverside"
$ws.Range("A24").Value = "UCSD"
$ws.Range("B24").Value = "University of California-San Diego"
$ws.Range("A25").Value = "UCSF"
$ws.Range("B25").Value = "University of California-San Francisco"
$ws.Range("A26").Value = "UCSD"
$ws.Range("B26").Value = "University of California-Santa Barbara"
$ws.Range("A27").Value = "UCSC"
$ws.Range("B27").Value = "University of California-Santa Cruz"
$ws.Range("A28").Value = "University of Colorado Denver"
$ws.Range("A29").Value = "UIUC"
$ws.Range("B29").Value = "University of Illinois Urbana-Champaign"
$ws.Range("A30").Value = "UMBC"
$ws.Range("B30").Value = "University of Maryland-Baltimore County"
$ws.Range("A31").Value = "UMCP"
$ws.Range("B31").Value = "University of Maryland-College Park"
$ws.Range("A32").Value = "UMass"
$ws.Range("B32").Value = "University of Massachusetts-Amherst"
$ws.Range("A33").Value = "UMass"
$ws.Range("B33").Value = "University of Massachusetts-Boston"
$ws.Range("A34").Value = "UMass"
$ws.Range("B34").Value = "University of Massachusetts-Lowell"
$ws.Range("A35").Value = "University of Michigan"
$ws.Range("A36").Value = "U MN Twin Cities"
$ws.Range("B36").Value = "University of Minnesota-Twin Cities"
$ws.Range("A37").Value = "U MO"
$ws.Range("B37").Value = "University of Missouri-Columbia"
$ws.Range("A38").Value = "U MO"
$ws.Range("B38").Value = "University of Missouri-Columbia"
$ws.Range("A39").Value = "U NE Omaha"
$ws.Range("B39").Value = "University of Nebraska Medical Center"
$ws.Range("A40").Value = "U NE Lincoln"
$ws.Range("B40").Value = "University of Nebraska-Lincoln"
$ws.Range("A41").Value = "U NV"
$ws.Range("B41").Value = "University of Nevada-Las Vegas"
$ws.Range("A42").Value = "U NV"
$ws.Range("B42").Value = "University of Nevada-Reno"
$ws.Range("A43").Value = "University of New Hampshire"
$ws.Range("A44").Value = "University of New Mexico"
$ws.Range("A45").Value = "U OK Norman"
$ws.Range("B45").Value = "University of Oklahoma"
$ws.Range("A46").Value = "University of Pittsburgh"
$ws.Range("A47").Value = "University of South Carolina"
$ws.Range("A48").Value = "University of Virginia"
$ws.Range("A49").Value = "University of Washington"
$ws.Range("A50").Value = "U WI Madison"
$ws.Range("B50").Value = "University of Wisconsin-Madison"
$ws.Range("A51").Value = "U WI Milwaukee"
$ws.Range("B51").Value = "University of Wisconsin-Milwaukee"
$ws.Range("A52").Value = "Cornell U"
$ws.Range("B52").Value = "Weill Medical College of Cornell University"

# --- Formatting: column A data rows (3-52) use plain (non-bold) Arial font ---
$ws.Range("A3:A52").Font.Name = "Arial"

# --- Formatting: stray empty styled cell C15 (kept blank, Arial-formatted) ---
$ws.Range("C15").Font.Name = "Arial"

# --- Remove now-stale trailing empty rows (previously 1000 & 1001) ---
$ws.Rows.Item(1000).Delete()
$ws.Rows.Item(1000).Delete()
